$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NATMI LR-pair (Sema4d-Plxnb1) TPM-derived stats for rows 2-10.
# Ligand stats (E-J) depend only on the sending cluster (rows 2-4 = ECs,
# 5-7 = FAPs, 8-10 = MuSCs) and receptor stats (K-P) only on the target
# cluster (rows 2/5/8 = ECs, 3/6/9 = FAPs, 4/7/10 = MuSCs); edge weights
# (Q-T) are recomputed per row from the new ligand/receptor values.

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2115286666666667
$ws.Range("H2").Value = 0.634586
$ws.Range("I2").Value = 0.08153347995807345
$ws.Range("J2").Value = 0.08153347995807345
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2807486666666667
$ws.Range("N2").Value = 0.8422460000000001
$ws.Range("O2").Value = 0.03211396410631209
$ws.Range("P2").Value = 0.03211396410631208
$ws.Range("Q2").Value = 0.05938639112844445
$ws.Range("R2").Value = 0.534477520156
$ws.Range("S2").Value = 0.002618363248836286
$ws.Range("T2").Value = 0.002618363248836286
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2115286666666667
$ws.Range("H3").Value = 0.634586
$ws.Range("I3").Value = 0.08153347995807345
$ws.Range("J3").Value = 0.08153347995807345
$ws.Range("O3").Value = 0.06996648921957034
$ws.Range("P3").Value = 0.06996648921957033
$ws.Range("Q3").Value = 0.1293847524062222
$ws.Range("R3").Value = 1.164462771656
$ws.Range("S3").Value = 0.0057046113465206
$ws.Range("T3").Value = 0.005704611346520599
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2115286666666667
$ws.Range("H4").Value = 0.634586
$ws.Range("I4").Value = 0.08153347995807345
$ws.Range("J4").Value = 0.08153347995807345
$ws.Range("O4").Value = 0.8979195466741177
$ws.Range("P4").Value = 0.8979195466741176
$ws.Range("Q4").Value = 1.660467739956889
$ws.Range("R4").Value = 14.944209659612
$ws.Range("S4").Value = 0.07321050536271656
$ws.Range("T4").Value = 0.07321050536271656
$ws.Range("I5").Value = 0.4735790235655714
$ws.Range("J5").Value = 0.4735790235655714
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2807486666666667
$ws.Range("N5").Value = 0.8422460000000001
$ws.Range("O5").Value = 0.03211396410631209
$ws.Range("P5").Value = 0.03211396410631208
$ws.Range("Q5").Value = 0.3449398840593333
$ws.Range("R5").Value = 3.104458956534
$ws.Range("S5").Value = 0.01520849976428709
$ws.Range("T5").Value = 0.01520849976428708
$ws.Range("I6").Value = 0.4735790235655714
$ws.Range("J6").Value = 0.4735790235655714
$ws.Range("O6").Value = 0.06996648921957034
$ws.Range("P6").Value = 0.06996648921957033
$ws.Range("S6").Value = 0.0331346616469152
$ws.Range("T6").Value = 0.03313466164691519
$ws.Range("I7").Value = 0.4735790235655714
$ws.Range("J7").Value = 0.4735790235655714
$ws.Range("O7").Value = 0.8979195466741177
$ws.Range("P7").Value = 0.8979195466741176
$ws.Range("S7").Value = 0.4252358621543691
$ws.Range("T7").Value = 0.4252358621543691
$ws.Range("I8").Value = 0.4448874964763552
$ws.Range("J8").Value = 0.4448874964763552
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2807486666666667
$ws.Range("N8").Value = 0.8422460000000001
$ws.Range("O8").Value = 0.03211396410631209
$ws.Range("P8").Value = 0.03211396410631208
$ws.Range("Q8").Value = 0.3240418891415556
$ws.Range("R8").Value = 2.916377002274
$ws.Range("S8").Value = 0.01428710109318872
$ws.Range("T8").Value = 0.01428710109318871
$ws.Range("I9").Value = 0.4448874964763552
$ws.Range("J9").Value = 0.4448874964763552
$ws.Range("O9").Value = 0.06996648921957034
$ws.Range("P9").Value = 0.06996648921957033
$ws.Range("S9").Value = 0.03112721622613454
$ws.Range("T9").Value = 0.03112721622613454
$ws.Range("I10").Value = 0.4448874964763552
$ws.Range("J10").Value = 0.4448874964763552
$ws.Range("O10").Value = 0.8979195466741177
$ws.Range("P10").Value = 0.8979195466741176
$ws.Range("S10").Value = 0.399473179157032
$ws.Range("T10").Value = 0.3994731791570319
